$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Arceneaux Strategies, LLC", 1),
    @("Arizona House of Representatives", 1),
    @("Bessemer Venture Partners", 1),
    @("Georgia Coordinated Campaign", 1),
    @("Presidential Inaugural Committee", 1),
    @("United States Senate, Office of Senator Robert P. Casey, Jr. (Retired)", 1)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
